# remove Gamelogic project, modify SLG building config
#
# The "Property" sheet lists BB_Build struct fields:
#   A=Id, B=Type, C=Public, D=Private, E=Save, F=View,
#   G=Index, H=SaveInterval, I=RelationValue, J=Desc
#
# Old row 7 held the "Desc" field (string / Friend / "描述").
# This edit inserts two new fields ahead of it:
#   row 7 -> "Icon"     (string / Friend / "图标")
#   row 8 -> "ShowName" (string / Friend / "名字")
# and pushes the original "Desc" field down to row 9 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BuildRow {
    param(
        [int]$Row,
        [string]$Id,
        [string]$Desc
    )

    $ws.Cells.Item($Row, 1).Value = $Id          # A: Id
    $ws.Cells.Item($Row, 2).Value = "string"     # B: Type
    $ws.Cells.Item($Row, 3).Value = $false        # C: Public
    $ws.Cells.Item($Row, 4).Value = $false        # D: Private
    $ws.Cells.Item($Row, 5).Value = $false        # E: Save
    $ws.Cells.Item($Row, 6).Value = $false        # F: View
    $ws.Cells.Item($Row, 7).Value = 0             # G: Index
    $ws.Cells.Item($Row, 8).Value = 0             # H: SaveInterval
    $ws.Cells.Item($Row, 9).Value = "Friend"      # I: RelationValue
    $ws.Cells.Item($Row, 10).Value = $Desc        # J: Desc

    # A, B, I and J use the sheet's "text" cell style (matches the rest
    # of the table, which stores ids/types/descriptions as plain text).
    $ws.Cells.Item($Row, 1).NumberFormat = "@"
    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 9).NumberFormat = "@"
    $ws.Cells.Item($Row, 10).NumberFormat = "@"
}

# Row 7 now becomes the new "Icon" field.
Set-BuildRow 7 "Icon" "图标"

# Row 8 is a brand-new "ShowName" field.
Set-BuildRow 8 "ShowName" "名字"

# Row 9 is the original "Desc" field, moved down from row 7.
Set-BuildRow 9 "Desc" "描述"

# The boolean "View" column drives a TRUE/FALSE list validation for every
# row below the table; that list now starts two rows later (F10 instead
# of F8) since two rows were inserted above it.
$ws.Range("F8:F1048576").Validation.Delete()
$ws.Range("F10:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Restore the author's last cell selection recorded in the saved file.
[void]$ws.Range("C13").Select()
